$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1 (row 7): "Wejdz na glowna strone sklepu..." - add cookie removal + screenshot step
$ws.Range("C7").Value = "Wejdź na główną stronę sklepu, usuń pliki cookies i zweryfikuj tytuł strony. Dodatkowo wykonaj zrzut ekranu."

# Expected / actual result for step 1 (row 7, columns E & F) - add cookies removed + screenshot saved
$ws.Range("E7").Value = "Strona główna sklepu wyświetlona`nTytuł strony: ''Koszulkifootball.sellingo.pl''`nPliki Cookies usunięte`nZrzut ekranu zapisany do pliku .jpg"
$ws.Range("F7").Value = "Strona główna sklepu wyświetlona`nTytuł strony: ''Koszulkifootball.sellingo.pl''`nPliki Cookies usunięte`nZrzut ekranu zapisany do pliku .jpg"

# Step 6 (row 12): hover-over menu tabs check - add screenshot instruction
$ws.Range("C12").Value = "Zweryfikuj czy każda z zakładek w menu, po najechaniu na nią myszką, zmienia kolor tekstu oraz tła na prawidowy. Dodatkowo wykonaj zrzut ekranu"

# Expected / actual result for step 6 (row 12, columns E & F) - add screenshot saved
$ws.Range("E12").Value = "Zakładki mają właściwości:`nkolor tekstu: rgba(255, 255, 255, 1)`nkolor tła: rgba(55, 55, 55, 1)`nZrzut ekranu zapisany do pliku .jpg"
$ws.Range("F12").Value = "Zakładki mają właściwości:`nkolor tekstu: rgba(255, 255, 255, 1)`nkolor tła: rgba(55, 55, 55, 1)`nZrzut ekranu zapisany do pliku .jpg"

# Update the view state to match the saved selection/scroll position
$ws.Range("D10").Select() | Out-Null
